$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 393 ("「笑顔」" post entry) - subsequent rows shift up
$ws.Rows.Item(393).Delete()
